$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cell values from the cryptos list refresh (price/volume
# updates plus a few rows that got reordered). Each new value is written
# with a leading apostrophe so Excel stores number-like strings (e.g.
# '8.38', '94.758.55') as literal text instead of auto-converting them
# to numbers/dates -- matching the original inlineStr text cells. The
# cell style is reset to Normal afterwards so no stray quote-prefix
# style gets attached.

$ws.Range('D2').Value = '''94.758.55'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  -3.45%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''3.436.10'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  +1.59%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '''  +0.08%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''238.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  -5.20%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''643.12'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  -2.28%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = '''  -0.94%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = '''0.406'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  -4.01%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = '''  +0.11%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''0.977'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  -5.75%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''3.434.92'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  +1.63%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = '''  -4.00%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''41.86'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  -2.14%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('E14').Value = '''  +3.19%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''94.631.70'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  -3.20%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''4.074.69'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  +1.48%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = '''  -1.21%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''8.38'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  -8.45%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''3.436.55'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  +2.38%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''17.58'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  -2.10%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''11.60'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  +4.84%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''0.490'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  -5.53%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''502.06'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  -1.43%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''3.25'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  -4.59%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''0.0000193'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  -3.55%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''6.55'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  -4.61%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''94.44'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  -2.12%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('B28').Value = '''WrappedeETH'
$ws.Range('B28').Style = 'Normal'
$ws.Range('C28').Value = '''https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('C28').Style = 'Normal'
$ws.Range('D28').Value = '''3.623.38'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  +1.57%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('B29').Value = '''Aptos'
$ws.Range('B29').Style = 'Normal'
$ws.Range('C29').Value = '''https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('C29').Style = 'Normal'
$ws.Range('D29').Value = '''11.98'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  -2.54%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''11.70'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  +0.72%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = '''  +0.14%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = '''2.76'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  +6.23%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = '''  -3.37%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('B34').Value = '''Cronos'
$ws.Range('B34').Style = 'Normal'
$ws.Range('C34').Value = '''https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('C34').Style = 'Normal'
$ws.Range('D34').Value = '''0.180'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  -4.34%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('B35').Value = '''Binance-PegBSC-USD'
$ws.Range('B35').Style = 'Normal'
$ws.Range('C35').Value = '''https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('C35').Style = 'Normal'
$ws.Range('D35').Value = '''0.998'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  -0.26%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''29.76'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  +3.63%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''0.554'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  -0.68%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''551.28'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  +4.56%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''7.72'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  -2.11%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''1.45'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  -1.10%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = '''  -0.21%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = '''  +0.06%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''0.916'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  +8.52%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''24.09'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  -1.33%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = '''  +0.15%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = '''3.70'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  +0.18%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = '''5.66'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  +2.49%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').Value = '''Stacks'
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Value = '''https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Value = '''2.22'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  -0.13%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Value = '''dogwifhat'
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = '''https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = '''3.34'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  +3.98%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('B50').Value = '''VeChain'
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = '''https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = '''0.0411'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  -3.12%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''55.00'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  -0.76%  '
$ws.Range('E51').Style = 'Normal'
